$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hlookup")

# ---------------------------------------------------------------
# 0. Fix up the formatting of a couple of existing cells so that
#    they match the rest of the (now larger) table. Do this first
#    so that later format-copies from these cells pick up the
#    corrected style.
# ---------------------------------------------------------------
$ws.Range("C7").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 2000000

$ws.Range("C7").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null

$ws.Range("C7").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null

$ws.Range("C7").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# 1. Extend the lookup table (rows 15-17) with three more rows
#    (18, 19, 20) holding data for city codes C02, C03 and C05.
#    Copy formatting from existing rows so the new cells pick up
#    the same styles that are already used in the sheet.
# ---------------------------------------------------------------

# Row 18 (C02 / Delhi / Delhi / 19000000)
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null

$ws.Range("C18").Value = "C02"
$ws.Range("D18").Value = 19000000
$ws.Range("E18").Value = "Delhi"
$ws.Range("F18").Value = "Delhi"

# Row 19 (C03 / karnataka / Bengaluru / 12000000)
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null

$ws.Range("C19").Value = "C03"
$ws.Range("D19").Value = 12000000
$ws.Range("E19").Value = "karnataka"
$ws.Range("F19").Value = "Bengaluru"

# Row 20 (C05 / Kolkata / west Bengal / 11000000)
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F20").PasteSpecial(-4122) | Out-Null

$ws.Range("C20").Value = "C05"
$ws.Range("D20").Value = 11000000
$ws.Range("E20").Value = "Kolkata"
$ws.Range("F20").Value = "west Bengal"

# ---------------------------------------------------------------
# 2. Rewrite the HLOOKUP formulas in rows 7-11 so that they pull
#    their results from the now-expanded lookup table.
# ---------------------------------------------------------------

# Row 7 - only D7 reference style changes (absolute -> mixed)
$ws.Range("D7").Formula = "=HLOOKUP(D6,`$C15:`$F17,2,FALSE)"

# Row 8 - was plain values, now formulas against the bigger table
$ws.Range("D8").Formula = "=HLOOKUP(D6,`$C`$15:`$F`$20,4,FALSE)"
$ws.Range("E8").Formula = "=HLOOKUP(E6,`$C`$15:`$F`$20,4,FALSE)"
$ws.Range("F8").Formula = "=HLOOKUP(F6,`$C`$15:`$F`$20,4,FALSE)"

# Row 9 - was plain values, now formulas against the bigger table
$ws.Range("D9").Formula = "=HLOOKUP(D6,`$C15:`$F20,5,FALSE)"
$ws.Range("E9").Formula = "=HLOOKUP(E6,`$C15:`$F20,5,FALSE)"
$ws.Range("F9").Formula = "=HLOOKUP(F6,`$C15:`$F20,5,FALSE)"

# Row 10 - formulas unchanged (still against the original C15:F17 table)
$ws.Range("D10").Formula = "=HLOOKUP(D6,`$C`$15:`$F`$17,3,FALSE)"
$ws.Range("E10").Formula = "=HLOOKUP(E6,`$C`$15:`$F`$17,3,FALSE)"
$ws.Range("F10").Formula = "=HLOOKUP(F6,`$C`$15:`$F`$17,3,FALSE)"

# Row 11 - was plain values, now formulas (note the typo'd ranges,
# preserved exactly as authored)
$ws.Range("D11").Formula = "=HLOOKUP(D6,`$C15:F420,6,FALSE)"
$ws.Range("E11").Formula = "=HLOOKUP(E6,`$C15:G420,6,FALSE)"
$ws.Range("F11").Formula = "=HLOOKUP(F6,`$C15:H420,6,FALSE)"

$wb.Application.Calculate()

# ---------------------------------------------------------------
# 3. Update the active selection shown in the sheet view.
# ---------------------------------------------------------------
$ws.Range("D11:F11").Select() | Out-Null
